$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D33:D62 numeric values (66..95) ---
$dVals = @(66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item(33 + $i, 4).Value = $dVals[$i]
}

# --- G68:G96 running counter formulas (extends existing G column tally) ---
$ws.Range("G67").Copy()
$ws.Range("G68:G96").PasteSpecial(-4122)
$ws.Range("G68").Formula = "=G67+1"
$ws.Range("G69:G96").Formula = "=G68+1"

# Trailing formatted-but-empty G cells (style carried one row further than data)
$ws.Range("G67").Copy()
$ws.Range("G97:G98").PasteSpecial(-4122)
$ws.Range("G97:G98").ClearContents()

# --- H67:H96 new tooltip text (imported categories tooltips) ---
$ws.Cells.Item(67, 8).Value = 'Vaske og tørre sig på kroppen og kropsdele med anvendelse af vand og passende rensemidler f.eks. tage bad, brusebad, vaske hænder og fødder, ansigt og hår og tørre sig med håndklæde.'
$ws.Cells.Item(68, 8).Value = 'Planlægge og udføre toiletbesøg til udskillelse af affaldsprodukter (menstruation, urin og afføring) og efterfølgende rengøring.'
$ws.Cells.Item(69, 8).Value = 'Pleje af de dele af kroppen, som behøver anden pleje end vask og tørring f.eks. hud, ansigt, tænder, hår, negle og kønsdele'
$ws.Cells.Item(70, 8).Value = 'Udføre sammensatte handlinger i forbindelse med på- og afklædning, at tage fodbeklædning på og af i rækkefølge'
$ws.Cells.Item(71, 8).Value = 'Udføre sammensatte handlinger i forbindelse med indtagelse af føde, som er serveret for én, få maden op til munden og spise på en kulturelt accepteret måde, skære eller bryde maden i stykker, åbne flasker og dåser, anvende spiseredskaber, deltage i måltider og i festligheder.'
$ws.Cells.Item(72, 8).Value = 'Holde fast om en drik, tage drikken op til munden og drikke på en kulturelt accepteret måde, blande, omrøre og skænke drikke op, åbne flasker og dåser, bruge sugerør eller drikke af rindende vand fra en hane eller en kilde; amning'
$ws.Cells.Item(73, 8).Value = 'Sikre sit velvære, helbred og fysiske og psykiske velbefindende ved f.eks. at indtage varieret kost, have passende niveau af fysisk aktivitet, holde sig varm eller afkølet, undgå skader på helbredet, dyrke sikker sex inkl. anvendelse af kondomer, lade sig vaccinere og følge regelmæssige helbredsundersøgelser.'
$ws.Cells.Item(74, 8).Value = 'Indtagelse og bearbejdning af fødemidler og væsker gennem munden. Inkl.: at suge, tygge, bide og behandle maden i mundhulen, spytflåd, at synke, gylpe, spytte og kaste op; tilstande som dysfagi, aspiration af føde,'
$ws.Cells.Item(75, 8).Value = 'Udføre simple, komplekse og sammensatte handlinger til planlægning, styring og gennemførelse af dagligt tilbagevendende rutiner eller pligter som f.eks. at overholde tider og lægge planer for særlige aktiviteter i løbet af dagen.'
$ws.Cells.Item(76, 8).Value = 'Vælge, tilvejebringe og transportere varer, som er nødvendige i dagliglivet som f.eks. at vælge, købe, transportere og opbevare mad, drikke, tøj, rengøringsmidler, brændsel, husholdningsgenstande og værktøj; tilvejebringe brugsgenstande og service.'
$ws.Cells.Item(77, 8).Value = 'Planlægge, tilberede og servere enkle eller sammensatte måltider til sig selv og andre som f.eks. at sammensætte et måltid, udvælge appetitlig mad og drikke, fremskaffe ingredienser til tilberedning af måltider;'
$ws.Cells.Item(78, 8).Value = 'Holde hus ved at gøre rent, vaske tøj, bruge husholdningsmaskiner, opbevare mad og smide affald ud, f.eks. ved at feje, moppe, tørre borde, […]'
$ws.Cells.Item(79, 8).Value = 'Skifte kropsstilling og bevæge sig fra et sted til et andet som f.eks. at flytte sig fra en stol til liggende stilling og skifte til og fra knælende eller hugsiddende stilling'
$ws.Cells.Item(80, 8).Value = 'Flytte sig fra en overflade til en anden som f.eks. at glide hen ad en bænk eller bevæge sig fra seng til stol uden at ændre kroppens stilling'
$ws.Cells.Item(81, 8).Value = 'Løfte en genstand op og flytte noget fra et sted til et andet som f.eks. at løfte en kop eller bære et barn fra et rum til et andet'
$ws.Cells.Item(82, 8).Value = 'Bevæge sig til fods skridt for skridt på et underlag, således at den ene fod hele tiden hviler på underlaget, som når man slentrer, går forlæns, baglæns eller sidelæns.'
$ws.Cells.Item(83, 8).Value = 'Bevæge sig fra et sted til et andet på andre måder end ved at gå […].'
$ws.Cells.Item(84, 8).Value = 'Gang og færden i forskellige omgivelser som f.eks. at gå mellem rum i huset, inden for en bygning eller ned ad gaden.'
$ws.Cells.Item(85, 8).Value = 'Bruge transportmidler som passager til at færdes omkring som f.eks. at blive kørt i en bil eller køre med […] taxi, bus, tog, sporvogn, undergrundsbane, skib eller fly'
$ws.Cells.Item(86, 8).Value = 'Funktioner bestemmende for respiratorisk og kardiovaskulær kapacitet, som er nødvendig ved fysisk anstrengelse.'
$ws.Cells.Item(87, 8).Value = 'Kraften som opstår ved kontraktion af en muskel eller en muskelgruppe.'
$ws.Cells.Item(88, 8).Value = 'Udvikle basale og komplekse kompetencer i sammensatte handlinger eller opgaver med det formål at påbegynde og gennemføre erhvervelsen af en færdighed, som f.eks. håndtering af værktøj eller spil som skak'
$ws.Cells.Item(89, 8).Value = 'Løsning af spørgsmål eller situationer ved at identificere og analysere emner, udvikle muligheder og løsninger, evaluere mulige virkninger af løsninger og gennemføre en valgt løsning som f.eks. ved løsning af en uoverensstemmelse mellem to personer.'
$ws.Cells.Item(90, 8).Value = 'Anvende udstyr, teknikker og andre midler med kommunikationsformål som f.eks. at ringe til en ven.'
$ws.Cells.Item(91, 8).Value = 'Overordnede mentale funktioner bestemmende for kendskab til og konstatering af relationerne til en selv, til andre, til tid, sted og andre omgivelser.'
$ws.Cells.Item(92, 8).Value = 'Overordnede mentale funktioner af fysiologisk og psykologisk art, som får personen til at opnå tilfredsstillelse af specifikke behov og overordnede mål på en vedholdende måde.'
$ws.Cells.Item(93, 8).Value = 'Specifikke mentale funktioner bestemmende for registrering, lagring genkaldelse af information efter behov.'
$ws.Cells.Item(94, 8).Value = 'Specifikke mentale funktioner forbundet med følelser og affektive komponenter i sindet.'
$ws.Cells.Item(95, 8).Value = 'Specifikke mentale funktioner først og fremmest knyttet til hjernens pandelapper omfattende kompleks målrettet adfærd som beslutningstagning, abstrakt tænkning, planlægning og gennemførelse af planer, mental fleksibilitet og tilpasning af adfærden efter omstændighederne, såkaldte eksekutive funktioner.'
$ws.Cells.Item(96, 8).Value = 'Deltage i alle aspekter af et arbejde, erhverv eller anden form for beskæftigelse […].'

# Match existing plain (column-default-styled) H cell formatting so the
# newly created H67:H96 cells do not carry a redundant explicit style index
$ws.Range("H60").Copy()
$ws.Range("H67:H96").PasteSpecial(-4122)

# --- Update sheet selection to match the authored view ---
$ws.Range("D77").Select()
